# Update the dSF column (F) values for the rows that were re-pulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = -2
$ws.Range("F15").Value = -2
$ws.Range("F17").Value = 4
$ws.Range("F22").Value = -5
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = -4
$ws.Range("F28").Value = -4
$ws.Range("F29").Value = 2
$ws.Range("F35").Value = -3
